# CreateUserData.xlsx edit: swap the sample user row from
#   name=Anuj Kumar / email=anuj2938@gmail.com / password=kishore@123
# to
#   name=virat / email=virat1989@gmail.com / password=virat99k
# keep both hyperlinks pointing at their original mailto targets, add a
# display tooltip (the old password text) to the password cell's hyperlink,
# and move the active selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rebuild the two hyperlinks (B2, C2) so they keep their original
# mailto targets/order (B2 -> rId1, C2 -> rId2) but C2 now also carries the
# old password text as its "display" attribute.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:anuj2938@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:kishore@123", $null, $null, "kishore@123")

# Adding a hyperlink resets the cell style - put the built-in "Hyperlink"
# style back on both cells so they keep looking like the original.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"

# --- Overwrite the row 2 values (do this after the hyperlink rebuild so the
# new text - not the stale snapshot the hyperlink edit captured - sticks).
$ws.Range("A2").Value = "virat"
$ws.Range("B2").Value = "virat1989@gmail.com"
$ws.Range("C2").Value = "virat99k"

# --- Move the selection to C3, matching the saved view state.
$ws.Range("C3").Select()
